$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Data")

# New column J: "ImagePath" header + sample image path value
$ws.Range("J1").Value = "ImagePath"
$ws.Range("J2").Value = "C:/Users/Public/Pictures/Sample Pictures/Penguins.jpg"

# Match the column's on-disk width (OOXML stores width = ColumnWidth + 5/6)
$ws.Columns.Item(10).ColumnWidth = (289/6)

# Scroll the view a bit and leave the new cell selected, like the authored session
$ws.Range("J2").Select()
$excel.ActiveWindow.ScrollColumn = 3
